$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date from 2023-10-05 (45204) to
# 2023-10-06 (45205) for every data row (rows 2 through 440).
$ws.Range("C2:C440").Value = (Get-Date -Year 2023 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0).Date
